# Update the "Coupling Parameters" sheet: switch data over from the DE
# case to the NL case, and adjust the related dynamic parameters.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coupling Parameters")

$ws.Range("B1").Value = "NL"
$ws.Range("B2").Value = 2019
$ws.Range("B4").Value = 2030
$ws.Range("B13").Value = 300

# Move the active selection to B3, as in the edited workbook.
$ws.Activate() | Out-Null
$ws.Range("B3").Select() | Out-Null
